$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "88.582.65"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3
$ws.Range("D3").Value = "3.274.98"
$ws.Range("E3").Value = "  -1.95%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "212.59"
$ws.Range("E5").Value = "  -3.07%  "

# Row 6
$ws.Range("D6").Value = "627.25"
$ws.Range("E6").Value = "  -1.71%  "

# Row 7
$ws.Range("E7").Value = "  +14.11%  "

# Row 8
$ws.Range("D8").Value = "0.718"
$ws.Range("E8").Value = "  +16.11%  "

# Row 9
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.08%  "

# Row 10
$ws.Range("D10").Value = "3.271.52"
$ws.Range("E10").Value = "  -2.35%  "

# Row 11
$ws.Range("D11").Value = "0.578"
$ws.Range("E11").Value = "  -5.66%  "

# Row 12
$ws.Range("E12").Value = "  +11.55%  "

# Row 13
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -5.21%  "

# Row 14
$ws.Range("D14").Value = "5.48"
$ws.Range("E14").Value = "  +0.91%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "34.13"
$ws.Range("E15").Value = "  -1.01%  "

# Row 16
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.877.76"
$ws.Range("E16").Value = "  -2.07%  "

# Row 17
$ws.Range("D17").Value = "88.548.52"
$ws.Range("E17").Value = "  +0.83%  "

# Row 18
$ws.Range("D18").Value = "3.316.68"
$ws.Range("E18").Value = "  -0.97%  "

# Row 19
$ws.Range("D19").Value = "3.17"
$ws.Range("E19").Value = "  -1.21%  "

# Row 20
$ws.Range("D20").Value = "14.06"
$ws.Range("E20").Value = "  -3.97%  "

# Row 21
$ws.Range("E21").Value = "  -3.03%  "

# Row 22
$ws.Range("D22").Value = "8.89"
$ws.Range("E22").Value = "  -2.71%  "

# Row 23
$ws.Range("D23").Value = "5.34"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").Value = "7.43"
$ws.Range("E24").Value = "  +0.36%  "

# Row 25
$ws.Range("D25").Value = "5.24"
$ws.Range("E25").Value = "  -3.16%  "

# Row 26
$ws.Range("D26").Value = "12.27"
$ws.Range("E26").Value = "  -0.59%  "

# Row 27
$ws.Range("D27").Value = "3.461.94"
$ws.Range("E27").Value = "  -1.67%  "

# Row 28
$ws.Range("D28").Value = "76.95"
$ws.Range("E28").Value = "  -2.52%  "

# Row 29
$ws.Range("D29").Value = "0.0000136"
$ws.Range("E29").Value = "  +3.91%  "

# Row 30
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.23%  "

# Row 31
$ws.Range("E31").Value = "  -5.15%  "

# Row 32
$ws.Range("E32").Value = "  +0.28%  "

# Row 33
$ws.Range("D33").Value = "563.85"
$ws.Range("E33").Value = "  -6.77%  "

# Row 34
$ws.Range("D34").Value = "8.76"
$ws.Range("E34").Value = "  -6.09%  "

# Row 35
$ws.Range("D35").Value = "1.37"
$ws.Range("E35").Value = "  -11.86%  "

# Row 36
$ws.Range("D36").Value = "7.15"
$ws.Range("E36").Value = "  +5.82%  "

# Row 37
$ws.Range("E37").Value = "  -4.68%  "

# Row 38
$ws.Range("D38").Value = "0.141"
$ws.Range("E38").Value = "  -7.60%  "

# Row 39
$ws.Range("D39").Value = "22.73"
$ws.Range("E39").Value = "  -3.21%  "

# Row 40
$ws.Range("D40").Value = "21.84"
$ws.Range("E40").Value = "  +2.07%  "

# Row 41
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.11%  "

# Row 42
$ws.Range("D42").Value = "3.11"
$ws.Range("E42").Value = "  +0.52%  "

# Row 43
$ws.Range("E43").Value = "  -4.71%  "

# Row 44
$ws.Range("D44").Value = "2.03"
$ws.Range("E44").Value = "  -2.00%  "

# Row 45
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$ws.Range("D46").Value = "153.81"
$ws.Range("E46").Value = "  -2.81%  "

# Row 47
$ws.Range("D47").Value = "180.90"
$ws.Range("E47").Value = "  -4.95%  "

# Row 48
$ws.Range("D48").Value = "0.135"
$ws.Range("E48").Value = "  +19.20%  "

# Row 49
$ws.Range("D49").Value = "44.78"
$ws.Range("E49").Value = "  -2.98%  "

# Row 50
$ws.Range("D50").Value = "1.32"
$ws.Range("E50").Value = "  -4.38%  "

# Row 51
$ws.Range("D51").Value = "4.22"
$ws.Range("E51").Value = "  -2.04%  "

Write-Host "Update complete"